$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.507.29"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").Value = "1.572.95"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  -1.61%  "

$ws.Range("D5").Value = "'211.22"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("E7").Value = "  -1.73%  "

$ws.Range("D8").Value = "'22.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.22%  "

$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("D11").Value = "'0.0870"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "1.797.38"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").Value = "1.560.51"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "27.489.09"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").Value = "'62.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.59%  "

$ws.Range("D18").Value = "'226.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.84%  "

$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("E23").Value = "  +2.33%  "

$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "'150.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.58%  "

$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("E28").Value = "  +1.44%  "

$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D33").Value = "1.458.93"
$ws.Range("E33").Value = "  +2.10%  "

$ws.Range("D34").Value = "'3.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.66%  "

$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("E42").Value = "  -3.13%  "

$ws.Range("E43").Value = "  -1.73%  "

$ws.Range("E44").Value = "  +7.15%  "

$ws.Range("D45").Value = "'0.974"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.06%  "

$ws.Range("D46").Value = "'64.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").Value = "1.709.26"
$ws.Range("E47").Value = "  +0.19%  "

$ws.Range("D48").Value = "'87.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("E49").Value = "  +3.12%  "

$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("E51").Value = "  -1.71%  "
